$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2034.8
$ws.Range("J17").Value = 2034.8
$ws.Range("L17").Value = 6104.4
$ws.Range("N17").Value = -6440.4
$ws.Range("H33").Value = 217.6
$ws.Range("I33").Value = 169.71428
$ws.Range("J33").Value = 888
$ws.Range("K33").Value = 169.71428
$ws.Range("L33").Value = 888
$ws.Range("M33").Value = 59.28572
$ws.Range("N33").Value = -1346
$ws.Range("H40").Value = 2855.125
$ws.Range("I40").Value = 1475
$ws.Range("J40").Value = 3315.1667
$ws.Range("K40").Value = 1475
$ws.Range("L40").Value = 3315.1667
$ws.Range("M40").Value = -1300
$ws.Range("N40").Value = -3665.1667
$ws.Range("H51").Value = 7572
$ws.Range("I51").Value = 9000
$ws.Range("J51").Value = 6960
$ws.Range("K51").Value = 9000
$ws.Range("L51").Value = 6960
$ws.Range("M51").Value = -8516
$ws.Range("N51").Value = -7928
$ws.Range("H53").Value = 3149
$ws.Range("I53").Value = 3149
$ws.Range("K53").Value = 3149
$ws.Range("M53").Value = -2512
$ws.Range("H88").Value = 1162
$ws.Range("I88").Value = 999
$ws.Range("K88").Value = 999
$ws.Range("M88").Value = -593
$ws.Range("H91").Value = 1162
$ws.Range("I91").Value = 999
$ws.Range("K91").Value = 999
$ws.Range("M91").Value = 405
$ws.Range("H112").Value = 4204
$ws.Range("J112").Value = 4333.2856
$ws.Range("L112").Value = 12999.8568
$ws.Range("N112").Value = -15215.8568
$ws.Range("H132").Value = 3742.1592
$ws.Range("I132").Value = 3689.6743
$ws.Range("K132").Value = 11069.0229
$ws.Range("M132").Value = -8539.0229
$ws.Range("H135").Value = 2741.68
$ws.Range("J135").Value = 3845.2
$ws.Range("L135").Value = 34606.8
$ws.Range("N135").Value = -39676.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1916.5538
$ws.Range("I32").Value = 1449.0159
$ws.Range("K32").Value = 1449.0159
$ws.Range("M32").Value = -1162.0159
$ws.Range("H63").Value = 2567.1177
$ws.Range("I63").Value = 2667.4614
$ws.Range("K63").Value = 2667.4614
$ws.Range("M63").Value = -1981.4614
$ws.Range("H66").Value = 2567.1177
$ws.Range("I66").Value = 2667.4614
$ws.Range("K66").Value = 13337.307
$ws.Range("M66").Value = -9905.307000000001
$ws.Range("H97").Value = 2461.261
$ws.Range("I97").Value = 1763.6154
$ws.Range("K97").Value = 1763.6154
$ws.Range("M97").Value = -1267.6154
$ws.Range("H132").Value = 1535.7894
$ws.Range("I132").Value = 1521.1111
$ws.Range("K132").Value = 4563.3333
$ws.Range("M132").Value = -2033.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4359
$ws.Range("J31").Value = 6832.8335
$ws.Range("L31").Value = 6832.8335
$ws.Range("N31").Value = -7422.8335
$ws.Range("H34").Value = 4359
$ws.Range("J34").Value = 6832.8335
$ws.Range("L34").Value = 6832.8335
$ws.Range("N34").Value = -7236.8335
$ws.Range("H50").Value = 40000
$ws.Range("J50").Value = 40000
$ws.Range("L50").Value = 40000
$ws.Range("N50").Value = -41250
$ws.Range("H62").Value = 19566.5
$ws.Range("I62").Value = 7998
$ws.Range("J62").Value = 31135
$ws.Range("K62").Value = 7998
$ws.Range("L62").Value = 31135
$ws.Range("M62").Value = -7374
$ws.Range("N62").Value = -32383
$ws.Range("H65").Value = 19566.5
$ws.Range("I65").Value = 7998
$ws.Range("J65").Value = 31135
$ws.Range("K65").Value = 39990
$ws.Range("L65").Value = 155675
$ws.Range("M65").Value = -36870
$ws.Range("N65").Value = -161915
$ws.Range("H86").Value = 6714.36
$ws.Range("I86").Value = 6031.5386
$ws.Range("J86").Value = 7454.0835
$ws.Range("K86").Value = 6031.5386
$ws.Range("L86").Value = 7454.0835
$ws.Range("M86").Value = -4908.5386
$ws.Range("N86").Value = -9700.083500000001
$ws.Range("H89").Value = 6714.36
$ws.Range("I89").Value = 6031.5386
$ws.Range("J89").Value = 7454.0835
$ws.Range("K89").Value = 30157.693
$ws.Range("L89").Value = 37270.4175
$ws.Range("M89").Value = -24541.693
$ws.Range("N89").Value = -48502.4175
$ws.Range("H132").Value = 3176.4
$ws.Range("I132").Value = 2153.4
$ws.Range("K132").Value = 6460.200000000001
$ws.Range("M132").Value = -3930.200000000001
$ws.Range("H141").Value = 244152.9
$ws.Range("J141").Value = 262668.2
$ws.Range("L141").Value = 262668.2
$ws.Range("N141").Value = -273028.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 696.1667
$ws.Range("I5").Value = 470
$ws.Range("J5").Value = 1148.5
$ws.Range("K5").Value = 1410
$ws.Range("L5").Value = 3445.5
$ws.Range("M5").Value = -1298
$ws.Range("N5").Value = -3669.5
$ws.Range("H37").Value = 139465.64
$ws.Range("J37").Value = 139465.64
$ws.Range("L37").Value = 418396.92
$ws.Range("N37").Value = -418620.92
$ws.Range("H131").Value = 1491.0303
$ws.Range("I131").Value = 965.38464
$ws.Range("J131").Value = 1832.7
$ws.Range("K131").Value = 2896.15392
$ws.Range("L131").Value = 5498.1
$ws.Range("M131").Value = 2143.84608
$ws.Range("N131").Value = -15578.1
$ws.Range("H135").Value = 696.1667
$ws.Range("I135").Value = 470
$ws.Range("J135").Value = 1148.5
$ws.Range("K135").Value = 4230
$ws.Range("L135").Value = 10336.5
$ws.Range("M135").Value = -1695
$ws.Range("N135").Value = -15406.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4251.25
$ws.Range("I80").Value = 4251.25
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 4251.25
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = -3253.25
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 4251.25
$ws.Range("I83").Value = 4251.25
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 21256.25
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = -16264.25
$ws.Range("M83").ClearContents()
$ws.Range("H107").Value = 432.75
$ws.Range("I107").Value = 239.3
$ws.Range("J107").Value = 1400
$ws.Range("K107").Value = 239.3
$ws.Range("L107").Value = 1400
$ws.Range("M107").Value = 1680.7
$ws.Range("N107").Value = -5240

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2022.3334
$ws.Range("I46").Value = 1267
$ws.Range("J46").Value = 2400
$ws.Range("K46").Value = 1267
$ws.Range("L46").Value = 2400
$ws.Range("M46").Value = -1079
$ws.Range("N46").Value = -2776
$ws.Range("H93").Value = 1153.8823
$ws.Range("I93").Value = 1093
$ws.Range("K93").Value = 1093
$ws.Range("M93").Value = 155

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8999
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 8999
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
